$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert a new blank column before column M ("mail") on the hidden
#    '#system' lookup sheet. This shifts columns M..AA (mail..xml) one
#    column to the right (becoming N..AB) and makes room for the new
#    "macro" command-category column at M.
# ---------------------------------------------------------------------------
$ws.Columns("M:M").Insert()

# 2. "target" column (A) - insert new category "macro" in alpha order
$ws.Cells.Item(1, 1).Value = "target"
$ws.Cells.Item(2, 1).Value = "aws.s3"
$ws.Cells.Item(3, 1).Value = "aws.ses"
$ws.Cells.Item(4, 1).Value = "base"
$ws.Cells.Item(5, 1).Value = "csv"
$ws.Cells.Item(6, 1).Value = "desktop"
$ws.Cells.Item(7, 1).Value = "excel"
$ws.Cells.Item(8, 1).Value = "external"
$ws.Cells.Item(9, 1).Value = "image"
$ws.Cells.Item(10, 1).Value = "io"
$ws.Cells.Item(11, 1).Value = "jms"
$ws.Cells.Item(12, 1).Value = "json"
$ws.Cells.Item(13, 1).Value = "macro"
$ws.Cells.Item(14, 1).Value = "mail"
$ws.Cells.Item(15, 1).Value = "number"
$ws.Cells.Item(16, 1).Value = "pdf"
$ws.Cells.Item(17, 1).Value = "rdbms"
$ws.Cells.Item(18, 1).Value = "redis"
$ws.Cells.Item(19, 1).Value = "sms"
$ws.Cells.Item(20, 1).Value = "sound"
$ws.Cells.Item(21, 1).Value = "ssh"
$ws.Cells.Item(22, 1).Value = "step"
$ws.Cells.Item(23, 1).Value = "web"
$ws.Cells.Item(24, 1).Value = "webalert"
$ws.Cells.Item(25, 1).Value = "webcookie"
$ws.Cells.Item(26, 1).Value = "ws"
$ws.Cells.Item(27, 1).Value = "ws.async"
$ws.Cells.Item(28, 1).Value = "xml"

# 3. "external" column (H) - fix typo runProgramNoWait, add new command
$ws.Cells.Item(1, 8).Value = "external"
$ws.Cells.Item(2, 8).Value = "runJUnit(className)"
$ws.Cells.Item(3, 8).Value = "runProgram(programPathAndParams)"
$ws.Cells.Item(4, 8).Value = "runProgramNoWait(programPathAndParams)"

# 4. "json" column (L) - insert beautify(json,var) / minify(json,var)
$ws.Cells.Item(1, 12).Value = "json"
$ws.Cells.Item(2, 12).Value = "addOrReplace(json,jsonpath,input,var)"
$ws.Cells.Item(3, 12).Value = "assertCorrectness(json,schema)"
$ws.Cells.Item(4, 12).Value = "assertElementCount(json,jsonpath,count)"
$ws.Cells.Item(5, 12).Value = "assertElementNotPresent(json,jsonpath)"
$ws.Cells.Item(6, 12).Value = "assertElementPresent(json,jsonpath)"
$ws.Cells.Item(7, 12).Value = "assertEqual(expected,actual)"
$ws.Cells.Item(8, 12).Value = "assertValue(json,jsonpath,expected)"
$ws.Cells.Item(9, 12).Value = "assertValues(json,jsonpath,array,exactOrder)"
$ws.Cells.Item(10, 12).Value = "assertWellformed(json)"
$ws.Cells.Item(11, 12).Value = "beautify(json,var)"
$ws.Cells.Item(12, 12).Value = "fromCsv(csv,header,jsonFile)"
$ws.Cells.Item(13, 12).Value = "minify(json,var)"
$ws.Cells.Item(14, 12).Value = "storeCount(json,jsonpath,var)"
$ws.Cells.Item(15, 12).Value = "storeValue(json,jsonpath,var)"
$ws.Cells.Item(16, 12).Value = "storeValues(json,jsonpath,var)"

# 5. "macro" column (M) - brand-new category, 3 commands
$ws.Cells.Item(1, 13).Value = "macro"
$ws.Cells.Item(2, 13).Value = "description()"
$ws.Cells.Item(3, 13).Value = "expects(var,default)"
$ws.Cells.Item(4, 13).Value = "produces(var,value)"

# 6. "web" column (W, was V before insert) - typo fix + 3 new commands
$ws.Cells.Item(1, 23).Value = "web"
$ws.Cells.Item(2, 23).Value = "assertAndClick(locator,label)"
$ws.Cells.Item(3, 23).Value = "assertAttribute(locator,attrName,value)"
$ws.Cells.Item(4, 23).Value = "assertAttributeContains(locator,attrName,contains)"
$ws.Cells.Item(5, 23).Value = "assertAttributeNotContains(locator,attrName,contains)"
$ws.Cells.Item(6, 23).Value = "assertAttributeNotPresent(locator,attrName)"
$ws.Cells.Item(7, 23).Value = "assertAttributePresent(locator,attrName)"
$ws.Cells.Item(8, 23).Value = "assertChecked(locator)"
$ws.Cells.Item(9, 23).Value = "assertContainCount(locator,text,count)"
$ws.Cells.Item(10, 23).Value = "assertCssNotPresent(locator,property)"
$ws.Cells.Item(11, 23).Value = "assertCssPresent(locator,property,value)"
$ws.Cells.Item(12, 23).Value = "assertElementByAttributes(nameValues)"
$ws.Cells.Item(13, 23).Value = "assertElementByText(locator,text)"
$ws.Cells.Item(14, 23).Value = "assertElementCount(locator,count)"
$ws.Cells.Item(15, 23).Value = "assertElementNotPresent(locator)"
$ws.Cells.Item(16, 23).Value = "assertElementPresent(locator)"
$ws.Cells.Item(17, 23).Value = "assertFocus(locator)"
$ws.Cells.Item(18, 23).Value = "assertFrameCount(count)"
$ws.Cells.Item(19, 23).Value = "assertFramePresent(frameName)"
$ws.Cells.Item(20, 23).Value = "assertIECompatMode()"
$ws.Cells.Item(21, 23).Value = "assertIENativeMode()"
$ws.Cells.Item(22, 23).Value = "assertLinkByLabel(label)"
$ws.Cells.Item(23, 23).Value = "assertNotChecked(locator)"
$ws.Cells.Item(24, 23).Value = "assertNotFocus(locator)"
$ws.Cells.Item(25, 23).Value = "assertNotText(locator,text)"
$ws.Cells.Item(26, 23).Value = "assertNotVisible(locator)"
$ws.Cells.Item(27, 23).Value = "assertOneMatch(locator)"
$ws.Cells.Item(28, 23).Value = "assertScrollbarHNotPresent(locator)"
$ws.Cells.Item(29, 23).Value = "assertScrollbarHPresent(locator)"
$ws.Cells.Item(30, 23).Value = "assertScrollbarVNotPresent(locator)"
$ws.Cells.Item(31, 23).Value = "assertScrollbarVPresent(locator)"
$ws.Cells.Item(32, 23).Value = "assertTable(locator,row,column,text)"
$ws.Cells.Item(33, 23).Value = "assertText(locator,text)"
$ws.Cells.Item(34, 23).Value = "assertTextContains(locator,text)"
$ws.Cells.Item(35, 23).Value = "assertTextCount(locator,text,count)"
$ws.Cells.Item(36, 23).Value = "assertTextList(locator,list,ignoreOrder)"
$ws.Cells.Item(37, 23).Value = "assertTextMatches(text,minMatch,scrollTo)"
$ws.Cells.Item(38, 23).Value = "assertTextNotPresent(text)"
$ws.Cells.Item(39, 23).Value = "assertTextOrder(locator,descending)"
$ws.Cells.Item(40, 23).Value = "assertTextPresent(text)"
$ws.Cells.Item(41, 23).Value = "assertTitle(text)"
$ws.Cells.Item(42, 23).Value = "assertValue(locator,value)"
$ws.Cells.Item(43, 23).Value = "assertValueOrder(locator,descending)"
$ws.Cells.Item(44, 23).Value = "assertVisible(locator)"
$ws.Cells.Item(45, 23).Value = "checkAll(locator)"
$ws.Cells.Item(46, 23).Value = "clearLocalStorage()"
$ws.Cells.Item(47, 23).Value = "click(locator)"
$ws.Cells.Item(48, 23).Value = "clickAndWait(locator,waitMs)"
$ws.Cells.Item(49, 23).Value = "clickByLabel(label)"
$ws.Cells.Item(50, 23).Value = "clickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(51, 23).Value = "clickOffset(locator,x,y)"
$ws.Cells.Item(52, 23).Value = "clickWithKeys(locator,keys)"
$ws.Cells.Item(53, 23).Value = "close()"
$ws.Cells.Item(54, 23).Value = "closeAll()"
$ws.Cells.Item(55, 23).Value = "deselect(locator,text)"
$ws.Cells.Item(56, 23).Value = "deselectMulti(locator,array)"
$ws.Cells.Item(57, 23).Value = "dismissInvalidCert()"
$ws.Cells.Item(58, 23).Value = "dismissInvalidCertPopup()"
$ws.Cells.Item(59, 23).Value = "doubleClick(locator)"
$ws.Cells.Item(60, 23).Value = "doubleClickAndWait(locator,waitMs)"
$ws.Cells.Item(61, 23).Value = "doubleClickByLabel(label)"
$ws.Cells.Item(62, 23).Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(63, 23).Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Cells.Item(64, 23).Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Cells.Item(65, 23).Value = "editLocalStorage(key,value)"
$ws.Cells.Item(66, 23).Value = "executeScript(var,script)"
$ws.Cells.Item(67, 23).Value = "focus(locator)"
$ws.Cells.Item(68, 23).Value = "goBack()"
$ws.Cells.Item(69, 23).Value = "goBackAndWait()"
$ws.Cells.Item(70, 23).Value = "maximizeWindow()"
$ws.Cells.Item(71, 23).Value = "mouseOver(locator)"
$ws.Cells.Item(72, 23).Value = "open(url)"
$ws.Cells.Item(73, 23).Value = "openAndWait(url,waitMs)"
$ws.Cells.Item(74, 23).Value = "openHttpBasic(url,username,password)"
$ws.Cells.Item(75, 23).Value = "openIgnoreTimeout(url)"
$ws.Cells.Item(76, 23).Value = "refresh()"
$ws.Cells.Item(77, 23).Value = "refreshAndWait()"
$ws.Cells.Item(78, 23).Value = "resizeWindow(width,height)"
$ws.Cells.Item(79, 23).Value = "saveAllWindowIds(var)"
$ws.Cells.Item(80, 23).Value = "saveAllWindowNames(var)"
$ws.Cells.Item(81, 23).Value = "saveAttribute(var,locator,attrName)"
$ws.Cells.Item(82, 23).Value = "saveAttributeList(var,locator,attrName)"
$ws.Cells.Item(83, 23).Value = "saveCount(var,locator)"
$ws.Cells.Item(84, 23).Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Cells.Item(85, 23).Value = "saveElement(var,locator)"
$ws.Cells.Item(86, 23).Value = "saveElements(var,locator)"
$ws.Cells.Item(87, 23).Value = "saveLocalStorage(var,key)"
$ws.Cells.Item(88, 23).Value = "saveLocation(var)"
$ws.Cells.Item(89, 23).Value = "savePageAs(var,sessionIdName,url)"
$ws.Cells.Item(90, 23).Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Cells.Item(91, 23).Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Cells.Item(92, 23).Value = "saveText(var,locator)"
$ws.Cells.Item(93, 23).Value = "saveTextArray(var,locator)"
$ws.Cells.Item(94, 23).Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Cells.Item(95, 23).Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Cells.Item(96, 23).Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Cells.Item(97, 23).Value = "saveValue(var,locator)"
$ws.Cells.Item(98, 23).Value = "scrollLeft(locator,pixel)"
$ws.Cells.Item(99, 23).Value = "scrollRight(locator,pixel)"
$ws.Cells.Item(100, 23).Value = "scrollTo(locator)"
$ws.Cells.Item(101, 23).Value = "select(locator,text)"
$ws.Cells.Item(102, 23).Value = "selectFrame(locator)"
$ws.Cells.Item(103, 23).Value = "selectMulti(locator,array)"
$ws.Cells.Item(104, 23).Value = "selectMultiOptions(locator)"
$ws.Cells.Item(105, 23).Value = "selectText(locator)"
$ws.Cells.Item(106, 23).Value = "selectWindow(winId)"
$ws.Cells.Item(107, 23).Value = "selectWindowAndWait(winId,waitMs)"
$ws.Cells.Item(108, 23).Value = "selectWindowByIndex(index)"
$ws.Cells.Item(109, 23).Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Cells.Item(110, 23).Value = "toggleSelections(locator)"
$ws.Cells.Item(111, 23).Value = "type(locator,value)"
$ws.Cells.Item(112, 23).Value = "typeKeys(locator,value)"
$ws.Cells.Item(113, 23).Value = "uncheckAll(locator)"
$ws.Cells.Item(114, 23).Value = "unselectAllText()"
$ws.Cells.Item(115, 23).Value = "upload(fieldLocator,file)"
$ws.Cells.Item(116, 23).Value = "verifyContainText(locator,text)"
$ws.Cells.Item(117, 23).Value = "verifyText(locator,text)"
$ws.Cells.Item(118, 23).Value = "wait(waitMs)"
$ws.Cells.Item(119, 23).Value = "waitForElementPresent(locator)"
$ws.Cells.Item(120, 23).Value = "waitForPopUp(winId,waitMs)"
$ws.Cells.Item(121, 23).Value = "waitForTextPresent(text)"
$ws.Cells.Item(122, 23).Value = "waitForTitle(text)"

# 7. "xml" column (AB, was AA before insert) - insert beautify/minify(xml,var)
$ws.Cells.Item(1, 28).Value = "xml"
$ws.Cells.Item(2, 28).Value = "assertCorrectness(xml,schema)"
$ws.Cells.Item(3, 28).Value = "assertElementCount(xml,xpath,count)"
$ws.Cells.Item(4, 28).Value = "assertElementNotPresent(xml,xpath)"
$ws.Cells.Item(5, 28).Value = "assertElementPresent(xml,xpath)"
$ws.Cells.Item(6, 28).Value = "assertValue(xml,xpath,expected)"
$ws.Cells.Item(7, 28).Value = "assertValues(xml,xpath,array,exactOrder)"
$ws.Cells.Item(8, 28).Value = "assertWellformed(xml)"
$ws.Cells.Item(9, 28).Value = "beautify(xml,var)"
$ws.Cells.Item(10, 28).Value = "minify(xml,var)"
$ws.Cells.Item(11, 28).Value = "storeCount(xml,xpath,var)"
$ws.Cells.Item(12, 28).Value = "storeValue(xml,xpath,var)"
$ws.Cells.Item(13, 28).Value = "storeValues(xml,xpath,var)"

# ---------------------------------------------------------------------------
# 8. Fix up the workbook-level defined names so each command-category name
#    still points at the right (now shifted / resized) range on '#system'.
# ---------------------------------------------------------------------------
$wb.Names.Item("external").RefersTo = "='#system'!`$H`$2:`$H`$4"
$wb.Names.Item("json").RefersTo = "='#system'!`$L`$2:`$L`$16"
$wb.Names.Item("mail").RefersTo = "='#system'!`$N`$2:`$N`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$O`$2:`$O`$15"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$P`$2:`$P`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$Q`$2:`$Q`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$R`$2:`$R`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$S`$2:`$S`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$T`$2:`$T`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$U`$2:`$U`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$V`$2:`$V`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$28"
$wb.Names.Item("web").RefersTo = "='#system'!`$W`$2:`$W`$122"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$Y`$2:`$Y`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$Z`$2:`$Z`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AB`$2:`$AB`$13"
$wb.Names.Add("macro", "='#system'!`$M`$2:`$M`$4")
